$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.803092333333333
$ws.Range("H2").Value = 5.409276999999999
$ws.Range("I2").Value = 0.1744886524959502
$ws.Range("J2").Value = 0.1744886524959502
$ws.Range("M2").Value = 6.045145666666667
$ws.Range("N2").Value = 18.135437
$ws.Range("O2").Value = 0.8160840232643366
$ws.Range("P2").Value = 0.8160840232643367
$ws.Range("Q2").Value = 10.89995580544989
$ws.Range("R2").Value = 98.09960224904901
$ws.Range("S2").Value = 0.1423974015428678
$ws.Range("T2").Value = 0.1423974015428678
$ws.Range("G3").Value = 1.803092333333333
$ws.Range("H3").Value = 5.409276999999999
$ws.Range("I3").Value = 0.1744886524959502
$ws.Range("J3").Value = 0.1744886524959502
$ws.Range("O3").Value = 0.09212864864242169
$ws.Range("P3").Value = 0.09212864864242169
$ws.Range("Q3").Value = 1.230508342267778
$ws.Range("R3").Value = 11.07457508041
$ws.Range("S3").Value = 0.01607540375788902
$ws.Range("T3").Value = 0.01607540375788901
$ws.Range("G4").Value = 1.803092333333333
$ws.Range("H4").Value = 5.409276999999999
$ws.Range("I4").Value = 0.1744886524959502
$ws.Range("J4").Value = 0.1744886524959502
$ws.Range("M4").Value = 0.6799149999999999
$ws.Range("N4").Value = 2.039745
$ws.Range("O4").Value = 0.09178732809324164
$ws.Range("P4").Value = 0.09178732809324165
$ws.Range("Q4").Value = 1.225949523818333
$ws.Range("R4").Value = 11.033545714365
$ws.Range("S4").Value = 0.01601584719519341
$ws.Range("T4").Value = 0.01601584719519341
$ws.Range("I5").Value = 0.4384883998568034
$ws.Range("J5").Value = 0.4384883998568034
$ws.Range("M5").Value = 6.045145666666667
$ws.Range("N5").Value = 18.135437
$ws.Range("O5").Value = 0.8160840232643366
$ws.Range("P5").Value = 0.8160840232643367
$ws.Range("Q5").Value = 27.39148999819645
$ws.Range("R5").Value = 246.523409983768
$ws.Range("S5").Value = 0.3578433775098813
$ws.Range("T5").Value = 0.3578433775098813
$ws.Range("I6").Value = 0.4384883998568034
$ws.Range("J6").Value = 0.4384883998568034
$ws.Range("O6").Value = 0.09212864864242169
$ws.Range("P6").Value = 0.09212864864242169
$ws.Range("S6").Value = 0.04039734372418515
$ws.Range("T6").Value = 0.04039734372418515
$ws.Range("I7").Value = 0.4384883998568034
$ws.Range("J7").Value = 0.4384883998568034
$ws.Range("M7").Value = 0.6799149999999999
$ws.Range("N7").Value = 2.039745
$ws.Range("O7").Value = 0.09178732809324164
$ws.Range("P7").Value = 0.09178732809324165
$ws.Range("Q7").Value = 3.080800025186667
$ws.Range("R7").Value = 27.72720022668
$ws.Range("S7").Value = 0.04024767862273695
$ws.Range("T7").Value = 0.04024767862273695
$ws.Range("G8").Value = 3.895605666666667
$ws.Range("H8").Value = 11.686817
$ws.Range("I8").Value = 0.3769851220961256
$ws.Range("J8").Value = 0.3769851220961256
$ws.Range("M8").Value = 6.045145666666667
$ws.Range("N8").Value = 18.135437
$ws.Range("O8").Value = 0.8160840232643366
$ws.Range("P8").Value = 0.8160840232643367
$ws.Range("Q8").Value = 23.54950371489211
$ws.Range("R8").Value = 211.9455334340291
$ws.Range("S8").Value = 0.3076515351510034
$ws.Range("T8").Value = 0.3076515351510034
$ws.Range("G9").Value = 3.895605666666667
$ws.Range("H9").Value = 11.686817
$ws.Range("I9").Value = 0.3769851220961256
$ws.Range("J9").Value = 0.3769851220961256
$ws.Range("O9").Value = 0.09212864864242169
$ws.Range("P9").Value = 0.09212864864242169
$ws.Range("Q9").Value = 2.658530116512223
$ws.Range("R9").Value = 23.92677104861
$ws.Range("S9").Value = 0.0347311298570144
$ws.Range("T9").Value = 0.0347311298570144
$ws.Range("G10").Value = 3.895605666666667
$ws.Range("H10").Value = 11.686817
$ws.Range("I10").Value = 0.3769851220961256
$ws.Range("J10").Value = 0.3769851220961256
$ws.Range("M10").Value = 0.6799149999999999
$ws.Range("N10").Value = 2.039745
$ws.Range("O10").Value = 0.09178732809324164
$ws.Range("P10").Value = 0.09178732809324165
$ws.Range("Q10").Value = 2.648680726851667
$ws.Range("R10").Value = 23.838126541665
$ws.Range("S10").Value = 0.03460245708810784
$ws.Range("T10").Value = 0.03460245708810784
$ws.Range("G11").Value = 0.1037266666666667
$ws.Range("H11").Value = 0.31118
$ws.Range("I11").Value = 0.01003782555112075
$ws.Range("J11").Value = 0.01003782555112075
$ws.Range("M11").Value = 6.045145666666667
$ws.Range("N11").Value = 18.135437
$ws.Range("O11").Value = 0.8160840232643366
$ws.Range("P11").Value = 0.8160840232643367
$ws.Range("Q11").Value = 0.6270428095177779
$ws.Range("R11").Value = 5.643385285660001
$ws.Range("S11").Value = 0.00819170906058418
$ws.Range("T11").Value = 0.00819170906058418
$ws.Range("G12").Value = 0.1037266666666667
$ws.Range("H12").Value = 0.31118
$ws.Range("I12").Value = 0.01003782555112075
$ws.Range("J12").Value = 0.01003782555112075
$ws.Range("O12").Value = 0.09212864864242169
$ws.Range("P12").Value = 0.09212864864242169
$ws.Range("Q12").Value = 0.07078757215555558
$ws.Range("R12").Value = 0.6370881494
$ws.Range("S12").Value = 0.0009247713033331266
$ws.Range("T12").Value = 0.0009247713033331264
$ws.Range("G13").Value = 0.1037266666666667
$ws.Range("H13").Value = 0.31118
$ws.Range("I13").Value = 0.01003782555112075
$ws.Range("J13").Value = 0.01003782555112075
$ws.Range("M13").Value = 0.6799149999999999
$ws.Range("N13").Value = 2.039745
$ws.Range("O13").Value = 0.09178732809324164
$ws.Range("P13").Value = 0.09178732809324165
$ws.Range("Q13").Value = 0.07052531656666666
$ws.Range("R13").Value = 0.6347278491
$ws.Range("S13").Value = 0.0009213451872034445
$ws.Range("T13").Value = 0.0009213451872034444
